# Weekly data refresh: a new week's observation is prepended to the data
# block (row 555, right after the fixed region that ends at row 554), which
# pushes every existing record down by one row. The oldest record (which was
# at row 669) ends up at the new last row, 670.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Insert a blank row at 555 - this shifts rows 555:669 down to 556:670 and
# grows the sheet's used range to row 670.
$ws.Rows.Item(555).Insert()

# The row that used to be 555 is now 556; pull the columns that never change
# between consecutive records in this block (identifiers/labels) from it so
# the new row 555 looks like a normal data row, then overwrite the columns
# that actually carry the new week's figures.
$ws.Range("A555").Value = $ws.Range("A556").Value()
$ws.Range("B555").Value = $ws.Range("B556").Value()
$ws.Range("C555").Value = $ws.Range("C556").Value()
$ws.Range("D555").Value = 45258
$ws.Range("E555").Value = $ws.Range("E556").Value()
$ws.Range("F555").Value = $ws.Range("F556").Value()
$ws.Range("G555").Value = $ws.Range("G556").Value()
$ws.Range("H555").Value = $ws.Range("H556").Value()
$ws.Range("I555").Value = $ws.Range("I556").Value()
$ws.Range("J555").Value = 510
$ws.Range("K555").Value = 9000
$ws.Range("L555").Value = 10000
$ws.Range("M555").Value = 9451
$ws.Range("N555").Value = $ws.Range("N556").Value()
$ws.Range("O555").Value = "Región Metropolitana"
$ws.Range("P555").Value = 189
$ws.Range("Q555").Value = $ws.Range("Q556").Value()
$ws.Range("R555").Value = $ws.Range("R556").Value()

"New dimension: $($ws.UsedRange.Rows.Count) rows"
